# Relatório 4 - spelling/proofing corrections and final review
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Davi Wei Tokikawa" - merge the 3 runs (split by spell-check
#    proofErr markers) into a single plain run.
# ---------------------------------------------------------------
$d.Content.Find.Execute("Davi Wei Tokikawa", $true, $false, $false, $false, $false, $true, 1, $false, "Davi Wei Tokikawa", 2) | Out-Null

# ---------------------------------------------------------------
# 2) "O pulso cosseno levantado..." paragraph - merge runs split by
#    proofErr markers (nyquist/tx/rx) into two runs.
#    First isolate the leading <w:tab/> into its own run so the
#    text replacement does not swallow it into the text content.
# ---------------------------------------------------------------
$p7 = $d.Paragraphs(7).Range
$splitTab7 = $d.Range($p7.Start + 1, $p7.Start + 1)
$d.Bookmarks.Add("TmpTabSplit7", $splitTab7) | Out-Null
$d.Bookmarks("TmpTabSplit7").Delete()

$d.Content.Find.Execute(
    "O pulso cosseno levantado é um pulso de nyquist, ou seja, tem a propriedade de que quando um símbolo é amostrado, os sinais dos outros símbolos ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O pulso cosseno levantado é um pulso de nyquist, ou seja, tem a propriedade de que quando um símbolo é amostrado, os sinais dos outros símbolos ",
    2) | Out-Null

$d.Content.Find.Execute(
    "vão estar cruzando o valor 0. Contudo, o emprego desse pulso é dado através do emprego de pulsos raiz cosseno levantado, de modo que, a transmissão (tx) gera esse pulso a partir de um filtro, transmite e a recepção (rx) reconstrói o sinal aplicando o mesmo tipo de filtro.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "vão estar cruzando o valor 0. Contudo, o emprego desse pulso é dado através do emprego de pulsos raiz cosseno levantado, de modo que, a transmissão (tx) gera esse pulso a partir de um filtro, transmite e a recepção (rx) reconstrói o sinal aplicando o mesmo tipo de filtro.",
    2) | Out-Null

# Re-split the boundary between the two merged runs (the two Find
# replacements above collapse back into a single run since both end
# up with identical, empty formatting).
$rngBoundary7 = $d.Content
$foundBoundary7 = $rngBoundary7.Find.Execute("dos outros símbolos ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundBoundary7) {
    $rngBoundary7.Collapse(0)
    $d.Bookmarks.Add("TmpBoundary7", $rngBoundary7) | Out-Null
    $d.Bookmarks("TmpBoundary7").Delete()
}

# ---------------------------------------------------------------
# 3) "O pulso raiz de cosseno levantado..." paragraph - merge runs
#    split by the "nyquist" proofErr marker into a single run.
# ---------------------------------------------------------------
$p8 = $d.Paragraphs(8).Range
$splitTab8 = $d.Range($p8.Start + 1, $p8.Start + 1)
$d.Bookmarks.Add("TmpTabSplit8", $splitTab8) | Out-Null
$d.Bookmarks("TmpTabSplit8").Delete()

$d.Content.Find.Execute(
    "O pulso raiz de cosseno levantado não é um pulso de nyquist, pois os sinais de outros símbolos não cruzam em zero quando um determinado símbolo corrente está sendo amostrado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O pulso raiz de cosseno levantado não é um pulso de nyquist, pois os sinais de outros símbolos não cruzam em zero quando um determinado símbolo corrente está sendo amostrado.",
    2) | Out-Null

# ---------------------------------------------------------------
# 4) "...respostas em frequência parecida e que conforme o aumento
#    do alfa, mas frequência..." -> "similares" / "mais", with a
#    _GoBack bookmark dropped right after "similares" (this is also
#    where the old trailing _GoBack bookmark effectively moves to).
# ---------------------------------------------------------------
$d.Content.Find.Execute("frequência parecida e que", $true, $false, $false, $false, $false, $true, 1, $false, "frequência similares e que", 2) | Out-Null

# Remove the old _GoBack bookmark (at the very end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Restore the original run boundary between " re" and "spostas em
# frequência " (they were two separate runs pre-edit and merged back
# together once their text content changed/matched formatting).
$rngReBoundary = $d.Content
$foundReBoundary = $rngReBoundary.Find.Execute("ambos têm re", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundReBoundary) {
    $rngReBoundary.Collapse(0)
    $d.Bookmarks.Add("TmpReBoundary", $rngReBoundary) | Out-Null
    $d.Bookmarks("TmpReBoundary").Delete()
}

# Split "spostas em frequência " from "similares".
$rngSimBoundary = $d.Content
$foundSimBoundary = $rngSimBoundary.Find.Execute("spostas em frequência ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSimBoundary) {
    $rngSimBoundary.Collapse(0)
    $d.Bookmarks.Add("TmpSimBoundary", $rngSimBoundary) | Out-Null
    $d.Bookmarks("TmpSimBoundary").Delete()
}

# Add the new _GoBack bookmark right after "similares".
$rngGoBack = $d.Content
$foundGoBack = $rngGoBack.Find.Execute("frequência similares", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundGoBack) {
    $rngGoBack.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rngGoBack) | Out-Null
}

# "mas" -> "mais": insert "i" between "ma" and "s", producing three
# runs: "...alfa, ma" | "i" | "s frequência banda..."
$rngMa = $d.Content
$foundMa = $rngMa.Find.Execute("alfa, ma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMa) {
    $rngMa.Collapse(0)
    $insertPoint = $rngMa.Start
    $rngMa.InsertAfter("i")

    $rngBeforeI = $d.Range($insertPoint, $insertPoint)
    $d.Bookmarks.Add("TmpBeforeI", $rngBeforeI) | Out-Null
    $d.Bookmarks("TmpBeforeI").Delete()

    $rngAfterI = $d.Range($insertPoint + 1, $insertPoint + 1)
    $d.Bookmarks.Add("TmpAfterI", $rngAfterI) | Out-Null
    $d.Bookmarks("TmpAfterI").Delete()
}
